$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 5): A5 = 100, B5 = 0.05
$ws.Range("A5").Value = 100
$ws.Range("B5").Value = 0.05

# Update the selection to D9 (matches the sheet view's active cell/sqref in the diff)
$ws.Range("D9").Select()
